$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 reuse the same header formatting as the
# existing header row (e.g. H1): bold font, thin border, centered.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-13 for the new columns I (I0) and J (IF)
$iValues = @(7, 1, 1, 1, 1, 1, 1, 13, 1, 7, 5, 4)
$jValues = @(7, 1, 5, 5, 5, 4, 3, 13, 2, 8, 6, 4)

for ($r = 0; $r -lt 12; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
